# Add two new employees ("Shyam" and "Ram") as rows 17 and 18 to both the
# "2025-02" and "2025-03" attendance sheets, mirroring the layout of the
# existing rows (Total Leaves = "0", Total WFH = "0", every daily
# attendance cell blank).

$wb = $excel.ActiveWorkbook

# Last data column differs per sheet: Feb sheet has 28 day columns (A:AE,
# 31 cols total), March sheet has 31 day columns (A:AH, 34 cols total).
$sheetLastCol = @{ "2025-02" = 31; "2025-03" = 34 }

# Row 17 -> Shyam, Row 18 -> Ram (appended after the last existing row, 16).
$newRows = @(
    @{ Row = 17; Name = "Shyam" },
    @{ Row = 18; Name = "Ram" }
)

foreach ($sheetName in $sheetLastCol.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $lastCol = $sheetLastCol[$sheetName]

    foreach ($entry in $newRows) {
        $r = $entry.Row
        $name = $entry.Name

        $ws.Cells.Item($r, 1).Value = $name

        # Leading apostrophe forces these numeric-looking strings to be
        # stored as text (matching "Total Leaves"/"Total WFH" = "0" text
        # cells used throughout the rest of the sheet).
        $ws.Cells.Item($r, 2).Value = "'0"
        $ws.Cells.Item($r, 3).Value = "'0"

        for ($c = 4; $c -le $lastCol; $c++) {
            # Leading apostrophe with nothing after it stores an empty
            # text cell, matching the blank daily-attendance cells.
            $ws.Cells.Item($r, $c).Value = "'"
        }
    }
}
